# Automatic update of files.
# Rotates the data of three observation records (rows 15, 16, 18) on the
# "Artfynd" sheet:
#   new row 15 <- old row 18
#   new row 16 <- old row 15
#   new row 18 <- old row 16
# Row 17 and all other rows are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 15 (becomes the old row 18 record) ----
$ws.Range("A15").Value = 130979082
$ws.Range("B15").Value = 57884
$ws.Range("E15").Value = 100109
$ws.Range("F15").Value = "Tretåig hackspett"
$ws.Range("G15").Value = "Picoides tridactylus"
$ws.Range("H15").Value = "(Linnaeus, 1758)"
$ws.Range("Q15").Value = 570952
$ws.Range("R15").Value = 6736563
$ws.Range("S15").Value = 1
$ws.Range("Z15").ClearContents()
$ws.Range("AB15").ClearContents()
$ws.Range("AC15").Value = "Äldre ringhack"
$ws.Range("AF15").ClearContents()
$ws.Range("AW15").Value = "Erik Danielsson"
$ws.Range("AX15").Value = "Erik Danielsson"

# ---- Row 16 (becomes the old row 15 record) ----
$ws.Range("A16").Value = 130983071
$ws.Range("P16").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q16").Value = 570817
$ws.Range("R16").Value = 6736417
$ws.Range("Z16").Value = "08:53"
$ws.Range("AB16").Value = "08:53"
$ws.Range("AF16").Value = ""
$ws.Range("AW16").Value = "Bo karlstens"
$ws.Range("AX16").Value = "Bo karlstens"

# ---- Row 18 (becomes the old row 16 record) ----
$ws.Range("A18").Value = 130983619
$ws.Range("B18").Value = 79244
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = "Garnlav"
$ws.Range("G18").Value = "Alectoria sarmentosa"
$ws.Range("H18").Value = "(Ach.) Ach."
$ws.Range("P18").Value = "Flytjärnsmyren, Dlr"
$ws.Range("Q18").Value = 570825
$ws.Range("R18").Value = 6736389
$ws.Range("S18").Value = 10
$ws.Range("Z18").Value = "08:54"
$ws.Range("AB18").Value = "08:54"
$ws.Range("AC18").ClearContents()
$ws.Range("AW18").Value = "Göran Ehn"
$ws.Range("AX18").Value = "Göran Ehn"
